# Clear the content of cell C1 (previously "~") while keeping its formatting,
# then select C1 as the active cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").ClearContents()
$ws.Range("C1").Select()
